$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Modelo" in F1
$ws.Range("F1").Value = "Modelo"

# Copy the formatting from the existing header cell (E1) onto F1 so the new
# header matches the other header cells (bold font, border, centered).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Add the model name value in F2 (plain, unstyled cell like the rest of row 2)
$ws.Range("F2").Value = "Pipeline(steps=[('model', GradientBoostingRegressor(n_estimators=150))])"
